# The list of "Requisitos" rows (B26:C29 / B26:C29) needs to be reordered so
# that the "LOM3231" requirement is listed first (right after "LOM3206" in the
# old order), i.e. it moves up from the 3rd position to the 1st position,
# shifting "LOM3206" and "LOM3215" down by one row each. "LOM3234" stays last.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lom3231 = "LOM3231 -  Métodos Experimentais da Física IV  (Indicação de Conjunto)`n"
$lom3206 = "LOM3206 -  Eletrônica  (Requisito)`n"
$lom3215 = "LOM3215 -  Física do Estado Sólido  (Requisito)`n"

$ws.Range("B26").Value = $lom3231
$ws.Range("C26").Value = $lom3231

$ws.Range("B27").Value = $lom3206
$ws.Range("C27").Value = $lom3206

$ws.Range("B28").Value = $lom3215
$ws.Range("C28").Value = $lom3215
